# Updated Ruby Examples to use common data folder
#
# - Refresh the cached "Update automatically" date field text (slide
#   master + every slide layout) from 2/12/2016 to 7/17/2014.
# - Move the "Table 3" graphic frame on slide 1 further down the slide.
# - Remove the stray "TEST TEST TEST" text box that was left on slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the cached datetimeFigureOut field text everywhere it shows
#    up: the slide master and each of its slide layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "7/17/2014"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2/3. Fix up slide 1: reposition the table and drop the leftover
#      "TextBox 4" shape containing "TEST TEST TEST".
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide.Shapes.Item($i)

    if ($shp.Name -eq "Table 3") {
        $shp.Left = 1524000 / 12700.0
        $shp.Top = 1397000 / 12700.0
    }
    elseif ($shp.Name -eq "TextBox 4") {
        $shp.Delete()
    }
}
